# Daily attendance processing - 2026-01-17 23:01:13
# Normalizes the "Recorded By" (column G) entries: whenever the
# comma-separated list of recorders ends with "System", the order of the
# list is reversed so "System" appears first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ",\s*"
        if ($parts[$parts.Length - 1] -eq "System") {
            $reversed = @()
            for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
